$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.063.93'
$ws.Range('D3').Value = '2.647.53'
$ws.Range('E3').Value = '  +1.49%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '531.71'
$ws.Range('E5').Value = '  +4.03%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '155.70'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.997'
$ws.Range('E8').Value = '  +0.85%  '
$ws.Range('E9').Value = '  -0.52%  '
$ws.Range('E10').Value = '  +4.78%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.352'
$ws.Range('E11').Value = '  +1.65%  '
$ws.Range('E12').Value = '  -0.06%  '
$ws.Range('D13').Value = '3.107.36'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').Value = '61.072.61'
$ws.Range('E14').Value = '  +1.19%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.03'
$ws.Range('E15').Value = '  +1.90%  '
$ws.Range('E16').Value = '  +2.57%  '
$ws.Range('D17').Value = '2.647.30'
$ws.Range('E17').Value = '  +1.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.77'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '355.07'
$ws.Range('E19').Value = '  +1.26%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.67'
$ws.Range('E20').Value = '  +0.79%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.24'
$ws.Range('E21').Value = '  +1.57%  '
$ws.Range('E22').Value = '  +0.29%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.79'
$ws.Range('E23').Value = '  +2.15%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.431'
$ws.Range('E24').Value = '  +2.11%  '
$ws.Range('E25').Value = '  +1.54%  '
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('D27').Value = '0.0₃0861'
$ws.Range('E27').Value = '  +2.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.40'
$ws.Range('E28').Value = '  +0.51%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.01%  '
$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.19'
$ws.Range('E30').Value = '  +7.48%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.63'
$ws.Range('E31').Value = '  +4.23%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.55'
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '150.10'
$ws.Range('E33').Value = '  -0.61%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.16'
$ws.Range('E34').Value = '  +3.89%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.21'
$ws.Range('E35').Value = '  +1.61%  '
$ws.Range('E36').Value = '  +9.23%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.898'
$ws.Range('E37').Value = '  +2.09%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '310.21'
$ws.Range('E38').Value = '  +5.21%  '
$ws.Range('E39').Value = '  +1.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.83'
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.648'
$ws.Range('E41').Value = '  +3.84%  '
$ws.Range('E42').Value = '  +1.57%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0564'
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('B45').Value = 'RenderToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.06'
$ws.Range('E45').Value = '  +3.78%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.91'
$ws.Range('E46').Value = '  +0.54%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0240'
$ws.Range('E47').Value = '  +2.64%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.33'
$ws.Range('E48').Value = '  +8.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.36'
$ws.Range('E49').Value = '  +0.34%  '
$ws.Range('D50').Value = '1.989.78'
$ws.Range('E50').Value = '  -0.56%  '
$ws.Range('E51').Value = '  +2.68%  '
